$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data set for rows 3..36 (columns A-D).
# Column B holds a fraction-of-day time value, except row 10 which stores
# the literal text "14:20" (kept as text by the source workbook).
$rows = @(
    @{ r = 3;  a = 42389; b = 0.5625;               c = 476;   d = 172 },
    @{ r = 4;  a = 42392; b = 0.45833333333333331;  c = 56;    d = 55 },
    @{ r = 5;  a = 42393; b = 0.40183238440735369;  c = 117;   d = 36 },
    @{ r = 6;  a = 42402; b = 0.625;                c = 498;   d = 74 },
    @{ r = 7;  a = 42403; b = 0.58368055555555554;  c = 430;   d = 76 },
    @{ r = 8;  a = 42404; b = 0.52777777777777779;  c = 256;   d = 208 },
    @{ r = 9;  a = 42436; b = 0.56944444444444442;  c = 541;   d = 998 },
    @{ r = 10; a = 42437; b = "14:20";               c = 551;   d = 825 },
    @{ r = 11; a = 42439; b = 0.65933814661409929;  c = 663;   d = 113 },
    @{ r = 12; a = 42441; b = 0.55815309248746492;  c = 4378;  d = 300 },
    @{ r = 13; a = 42442; b = 0.51066760248567022;  c = 3688;  d = 259 },
    @{ r = 14; a = 42443; b = 0.53119982396714727;  c = 3190;  d = 308 },
    @{ r = 15; a = 42444; b = 0.51742016319439443;  c = 2637;  d = 191 },
    @{ r = 16; a = 42446; b = 0.56646408071301035;  c = 2840;  d = 134 },
    @{ r = 17; a = 42466; b = 0.54166666666666663;  c = 473;   d = 23 },
    @{ r = 18; a = 42719; b = 0.69444444444444453;  c = 79;    d = 14 },
    @{ r = 19; a = 42720; b = 0.72222222222222221;  c = 233;   d = 194 },
    @{ r = 20; a = 42740; b = 0.61458333333333326;  c = 205;   d = 97 },
    @{ r = 21; a = 42744; b = 0.72934110256754703;  c = 12451; d = 1166 },
    @{ r = 22; a = 42745; b = 0.60417623205853732;  c = 3801;  d = 668 },
    @{ r = 23; a = 42746; b = 0.66666817603807427;  c = 13750; d = 665 },
    @{ r = 24; a = 42749; b = 0.64930555555555558;  c = 3891;  d = 215 },
    @{ r = 25; a = 42754; b = 0.52083333333333337;  c = 6730;  d = 319 },
    @{ r = 26; a = 42755; b = 0.44791666666666669;  c = 6730;  d = 621 },
    @{ r = 27; a = 42756; b = 0.3923611111111111;   c = 10109; d = 617 },
    @{ r = 28; a = 42758; b = 0.49652777777777779;  c = 9896;  d = 450 },
    @{ r = 29; a = 42761; b = 0.47934027777777777;  c = 3453;  d = 66 },
    @{ r = 30; a = 42767; b = 0.54185185925636103;  c = 3053;  d = 31 },
    @{ r = 31; a = 42774; b = 0.57639118479022766;  c = 11060; d = 549 },
    @{ r = 32; a = 42776; b = 0.61458485941221541;  c = 11505; d = 408 },
    @{ r = 33; a = 42787; b = 0.5590286977733987;   c = 13460; d = 422 },
    @{ r = 34; a = 42810; b = 0.64976728373393811;  c = 2942;  d = 20 },
    @{ r = 35; a = 42829; b = 0.65605302031566826;  c = 259;   d = 21 },
    @{ r = 36; a = 42851; b = 0.60451388888888891;  c = 223;   d = 17 }
)

foreach ($row in $rows) {
    $r = $row.r

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = $row.a
    $cellA.NumberFormat = "yyyy\-mm\-dd;@"
    $cellA.HorizontalAlignment = -4152

    $cellB = $ws.Cells.Item($r, 2)
    $cellB.Value = $row.b
    $cellB.NumberFormat = "h:mm;@"
    $cellB.HorizontalAlignment = -4152

    $ws.Cells.Item($r, 3).Value = $row.c
    $ws.Cells.Item($r, 4).Value = $row.d
}

$ws.Range("E36").Select()
